$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows above the current row 267 (old rows 267-286
# shift down to become rows 271-290).
$ws.Range("A267:T270").EntireRow.Insert()

# Common (constant) values shared by every data row in this block.
$mercadoId = 5
$mercado   = "Macroferia Regional de Talca"
$region    = "Maule"
$codreg    = 7
$tipo      = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$categoriaId = 100103004
$categoria   = "Durazno"
$origen      = "Región de O'Higgins"

# Row 267: Carson / Especial
$ws.Cells.Item(267, 1).Value = $mercadoId
$ws.Cells.Item(267, 2).Value = $mercado
$ws.Cells.Item(267, 3).Value = $region
$ws.Cells.Item(267, 4).Value = 44585
$ws.Cells.Item(267, 5).Value = $codreg
$ws.Cells.Item(267, 6).Value = $tipo
$ws.Cells.Item(267, 7).Value = $productoId
$ws.Cells.Item(267, 8).Value = $producto
$ws.Cells.Item(267, 9).Value = $categoriaId
$ws.Cells.Item(267, 10).Value = $categoria
$ws.Cells.Item(267, 11).Value = "Carson"
$ws.Cells.Item(267, 12).Value = "Especial"
$ws.Cells.Item(267, 13).Value = 400
$ws.Cells.Item(267, 14).Value = 12000
$ws.Cells.Item(267, 15).Value = 12000
$ws.Cells.Item(267, 16).Value = 12000
$ws.Cells.Item(267, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(267, 18).Value = $origen
$ws.Cells.Item(267, 19).Value = 667
$ws.Cells.Item(267, 20).Value = 18

# Row 268: Carson / Primera
$ws.Cells.Item(268, 1).Value = $mercadoId
$ws.Cells.Item(268, 2).Value = $mercado
$ws.Cells.Item(268, 3).Value = $region
$ws.Cells.Item(268, 4).Value = 44585
$ws.Cells.Item(268, 5).Value = $codreg
$ws.Cells.Item(268, 6).Value = $tipo
$ws.Cells.Item(268, 7).Value = $productoId
$ws.Cells.Item(268, 8).Value = $producto
$ws.Cells.Item(268, 9).Value = $categoriaId
$ws.Cells.Item(268, 10).Value = $categoria
$ws.Cells.Item(268, 11).Value = "Carson"
$ws.Cells.Item(268, 12).Value = "Primera"
$ws.Cells.Item(268, 13).Value = 400
$ws.Cells.Item(268, 14).Value = 10000
$ws.Cells.Item(268, 15).Value = 10000
$ws.Cells.Item(268, 16).Value = 10000
$ws.Cells.Item(268, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(268, 18).Value = $origen
$ws.Cells.Item(268, 19).Value = 556
$ws.Cells.Item(268, 20).Value = 18

# Row 269: Polar King / Especial
$ws.Cells.Item(269, 1).Value = $mercadoId
$ws.Cells.Item(269, 2).Value = $mercado
$ws.Cells.Item(269, 3).Value = $region
$ws.Cells.Item(269, 4).Value = 44585
$ws.Cells.Item(269, 5).Value = $codreg
$ws.Cells.Item(269, 6).Value = $tipo
$ws.Cells.Item(269, 7).Value = $productoId
$ws.Cells.Item(269, 8).Value = $producto
$ws.Cells.Item(269, 9).Value = $categoriaId
$ws.Cells.Item(269, 10).Value = $categoria
$ws.Cells.Item(269, 11).Value = "Polar King"
$ws.Cells.Item(269, 12).Value = "Especial"
$ws.Cells.Item(269, 13).Value = 300
$ws.Cells.Item(269, 14).Value = 12000
$ws.Cells.Item(269, 15).Value = 12000
$ws.Cells.Item(269, 16).Value = 12000
$ws.Cells.Item(269, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(269, 18).Value = $origen
$ws.Cells.Item(269, 19).Value = 800
$ws.Cells.Item(269, 20).Value = 15

# Row 270: Polar King / Primera
$ws.Cells.Item(270, 1).Value = $mercadoId
$ws.Cells.Item(270, 2).Value = $mercado
$ws.Cells.Item(270, 3).Value = $region
$ws.Cells.Item(270, 4).Value = 44585
$ws.Cells.Item(270, 5).Value = $codreg
$ws.Cells.Item(270, 6).Value = $tipo
$ws.Cells.Item(270, 7).Value = $productoId
$ws.Cells.Item(270, 8).Value = $producto
$ws.Cells.Item(270, 9).Value = $categoriaId
$ws.Cells.Item(270, 10).Value = $categoria
$ws.Cells.Item(270, 11).Value = "Polar King"
$ws.Cells.Item(270, 12).Value = "Primera"
$ws.Cells.Item(270, 13).Value = 400
$ws.Cells.Item(270, 14).Value = 10000
$ws.Cells.Item(270, 15).Value = 10000
$ws.Cells.Item(270, 16).Value = 10000
$ws.Cells.Item(270, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(270, 18).Value = $origen
$ws.Cells.Item(270, 19).Value = 667
$ws.Cells.Item(270, 20).Value = 15
